# Add a 28th row (all zeros) to the digit pattern matrix on Sheet1,
# matching the other 27 rows already present (A:AB), and update the
# active selection to reflect where the author left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill new row 28 across columns A:AB with 0, extending the pattern grid.
$ws.Range("A28:AB28").Value = 0

# Match the final selection recorded in the saved workbook.
$ws.Range("AE27").Select() | Out-Null
